$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update table title for November 2016
$ws.Range("A1").Value = "Table C.1 Average Heat Content of Fossil-Fuel Receipts, November 2016"

# Update data values (134 cell changes)
$ws.Range("B3").Value = 22.61
$ws.Range("C3").Value = 6.23
$ws.Range("C4").Value = 5.8
$ws.Range("B5").Value = 24.91
$ws.Range("C5").Value = 6.24
$ws.Range("B6").Value = 22.43
$ws.Range("C6").Value = 6.27
$ws.Range("B7").Value = "--"
$ws.Range("C7").Value = 5.74
$ws.Range("C8").Value = 5.79
$ws.Range("B10").Value = 24.61
$ws.Range("C10").Value = 6.16
$ws.Range("E10").Value = 1.04
$ws.Range("B11").Value = 26.33
$ws.Range("C11").Value = 5.67
$ws.Range("B12").Value = 25.93
$ws.Range("C12").Value = 6.33
$ws.Range("B13").Value = 24.48
$ws.Range("C13").Value = 5.78
$ws.Range("B14").Value = 20.13
$ws.Range("C14").Value = 5.81
$ws.Range("D14").Value = 27.42
$ws.Range("B15").Value = 17.68
$ws.Range("C15").Value = 5.8
$ws.Range("E15").Value = 1.01
$ws.Range("B16").Value = 22.5
$ws.Range("B17").Value = 18.5
$ws.Range("C17").Value = 5.89
$ws.Range("D17").Value = 26.52
$ws.Range("B18").Value = 24.63
$ws.Range("C18").Value = 5.8
$ws.Range("D18").Value = 27.83
$ws.Range("E18").Value = 1.06
$ws.Range("B19").Value = 18.05
$ws.Range("C19").Value = 5.85
$ws.Range("D19").Value = 27.13
$ws.Range("B20").Value = 16.63
$ws.Range("C20").Value = 5.81
$ws.Range("B21").Value = 17.579999999999998
$ws.Range("C21").Value = 5.79
$ws.Range("B22").Value = 17.149999999999999
$ws.Range("C22").Value = 5.78
$ws.Range("B23").Value = 17.579999999999998
$ws.Range("C23").Value = 5.8
$ws.Range("E23").Value = 1.06
$ws.Range("B24").Value = 17.55
$ws.Range("C24").Value = 5.8
$ws.Range("B25").Value = 16.89
$ws.Range("C25").Value = "--"
$ws.Range("E25").Value = 1.06
$ws.Range("B26").Value = 13.22
$ws.Range("C26").Value = 5.98
$ws.Range("B27").Value = 16.579999999999998
$ws.Range("C27").Value = 6
$ws.Range("B28").Value = 23.71
$ws.Range("C28").Value = 6.06
$ws.Range("D28").Value = 27.84
$ws.Range("B29").Value = 25.76
$ws.Range("C29").Value = 5.5
$ws.Range("B31").Value = 23.56
$ws.Range("C31").Value = 5.77
$ws.Range("D31").Value = 28.03
$ws.Range("B32").Value = 20.39
$ws.Range("C32").Value = 5.9
$ws.Range("D32").Value = 25.91
$ws.Range("B33").Value = 25.19
$ws.Range("B34").Value = 24.85
$ws.Range("C34").Value = 5.78
$ws.Range("E34").Value = 1.03
$ws.Range("B35").Value = 25.26
$ws.Range("C35").Value = 5.86
$ws.Range("B36").Value = 23.33
$ws.Range("C36").Value = 6.25
$ws.Range("E36").Value = 1.06
$ws.Range("B37").Value = 24.64
$ws.Range("C37").Value = 5.76
$ws.Range("E37").Value = 1.08
$ws.Range("B38").Value = 20.71
$ws.Range("D38").Value = 28.12
$ws.Range("B39").Value = 19.28
$ws.Range("C39").Value = 5.62
$ws.Range("B40").Value = 21.84
$ws.Range("C40").Value = 5.84
$ws.Range("D40").Value = 28.12
$ws.Range("E40").Value = 1.06
$ws.Range("B41").Value = 14.75
$ws.Range("C41").Value = 5.81
$ws.Range("E41").Value = 1.04
$ws.Range("B42").Value = 21.67
$ws.Range("B43").Value = 16.059999999999999
$ws.Range("C43").Value = 5.85
$ws.Range("D43").Value = 28.61
$ws.Range("B44").Value = 17.41
$ws.Range("C44").Value = 5.87
$ws.Range("B45").Value = 16.190000000000001
$ws.Range("C45").Value = "--"
$ws.Range("D45").Value = 28.61
$ws.Range("E45").Value = 1.03
$ws.Range("B46").Value = 17.34
$ws.Range("C46").Value = "--"
$ws.Range("B47").Value = 15.73
$ws.Range("C47").Value = 5.78
$ws.Range("E47").Value = 1.02
$ws.Range("B48").Value = 18.7
$ws.Range("B49").Value = 19.510000000000002
$ws.Range("C49").Value = 5.65
$ws.Range("C50").Value = "--"
$ws.Range("E50").Value = 1.08
$ws.Range("B52").Value = 16.829999999999998
$ws.Range("C52").Value = 5.92
$ws.Range("B53").Value = 20.47
$ws.Range("C53").Value = "--"
$ws.Range("B54").Value = 18.57
$ws.Range("B55").Value = 21.4
$ws.Range("C55").Value = 5.88
$ws.Range("B56").Value = 17.399999999999999
$ws.Range("C56").Value = 5.81
$ws.Range("E56").Value = 1.03
$ws.Range("B57").Value = 17.600000000000001
$ws.Range("C57").Value = 5.92
$ws.Range("E57").Value = 1.03
$ws.Range("B58").Value = 22.94
$ws.Range("B59").Value = 17.239999999999998
$ws.Range("B60").Value = 17.170000000000002
$ws.Range("C60").Value = 5.92
$ws.Range("E60").Value = 1.1000000000000001
$ws.Range("B61").Value = 18.95
$ws.Range("C61").Value = 6.14
$ws.Range("C62").Value = 5.6
$ws.Range("B63").Value = 19.47
$ws.Range("C63").Value = 6.14
$ws.Range("B64").Value = 19.190000000000001
$ws.Range("C64").Value = 6.08
$ws.Range("D64").Value = 28.14
